$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row103 = @(402, 253, 141, 4, 4, 102, 145, 62, 0)
$row104 = @(402, 253, 141, 4, 4, 103, 145, 0, 0)

for ($i = 0; $i -lt 9; $i++) {
    $ws.Cells.Item(103, $i + 1).Value = $row103[$i]
    $ws.Cells.Item(104, $i + 1).Value = $row104[$i]
}
